# Actualización de tablas de salida dir_cv*100 -SG
# Multiply the "dir_cv" columns (D = Dir_Educacion_cv, F = Dir_Empleo_cv,
# H = Dir_IPM_cv) by 100 for all data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(4, 6, 8)   # D, F, H
$firstRow = 2
$lastRow = 25

for ($r = $firstRow; $r -le $lastRow; $r++) {
    foreach ($c in $cols) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($null -ne $val) {
            $cell.Value = $val * 100
        }
    }
}
